$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "Done"
$ws.Range("E11").Value = "Throwing Knife"
$ws.Range("F11").Value = "Zidane Dagger"

$ws.Range("B12").Value = "InProgress"
$ws.Range("C12").Value = 22
$ws.Range("D12").Value = 23
$ws.Range("C12:D12").HorizontalAlignment = -4152
$ws.Range("E12").Value = "Axe"
$ws.Range("F12").Value = "Lani_ZidaneDagger"
$ws.Range("H12").Value = "Done"
$ws.Range("I12").Value = "Done"
